# Auto-generated edit: add data for 2024-09-03 to violent-crime-full-year workbook
$wb = $excel.ActiveWorkbook

$edits = @{
    'Citywide Totals' = @{
        "K2" = 5403
        "K3" = 5569
        "C4" = 1851
        "F4" = 1914
        "K4" = 1160
        "K5" = 399
        "K6" = 6184
        "C7" = 28396
        "F7" = 24107
        "K7" = 18715
    }
    'Logan Square' = @{
        "K6" = 105
        "K7" = 238
    }
    'Austin' = @{
        "K2" = 343
        "K3" = 376
        "K6" = 421
        "K7" = 1248
    }
    'South Chicago' = @{
        "K3" = 151
        "K6" = 95
        "K7" = 418
    }
    'Garfield Park' = @{
        "K3" = 295
        "K7" = 802
    }
    'West Pullman' = @{
        "K2" = 106
        "K5" = 14
        "K7" = 317
    }
    'Grand Crossing' = @{
        "K3" = 207
        "K6" = 185
        "K7" = 630
    }
    'New City' = @{
        "K2" = 137
        "K4" = 15
        "K7" = 426
    }
    'Fuller Park' = @{
        "K6" = 26
        "K7" = 74
    }
    'By Neighborhood' = @{
        "K2" = 164
        "K4" = 71
        "K8" = 1248
        "K11" = 356
        "K19" = 549
        "K20" = 434
        "K24" = 55
        "K27" = 176
        "K29" = 1005
        "K30" = 74
        "K33" = 802
        "K35" = 30
        "K36" = 247
        "K37" = 630
        "K42" = 693
        "K43" = 166
        "K44" = 158
        "K46" = 39
        "K47" = 128
        "K48" = 237
        "K52" = 490
        "K53" = 238
        "K54" = 363
        "K59" = 31
        "C63" = 280
        "F63" = 198
        "K63" = 51
        "K64" = 122
        "K65" = 426
        "K67" = 712
        "K73" = 163
        "K76" = 258
        "K78" = 217
        "K79" = 471
        "K83" = 418
        "K84" = 142
        "K85" = 885
        "K86" = 126
        "K87" = 31
        "K88" = 206
        "K89" = 272
        "K90" = 170
        "K91" = 208
        "K92" = 69
        "K94" = 253
        "K95" = 317
        "K98" = 87
        "C101" = 28396
        "F101" = 24107
        "K101" = 18715
    }
    'North Lawndale' = @{
        "K2" = 203
        "K3" = 256
        "K4" = 40
        "K7" = 712
    }
    'South Deering' = @{
        "K2" = 47
        "K7" = 142
    }
    'Lincoln Park' = @{
        "K2" = 22
        "K6" = 52
    }
    'Loop' = @{
        "K6" = 194
        "K7" = 363
    }
    'Englewood' = @{
        "K2" = 286
        "K3" = 361
        "K6" = 283
        "K7" = 1005
    }
    'Lake View' = @{
        "K6" = 116
        "K7" = 237
    }
    'Chatham' = @{
        "K2" = 163
        "K3" = 173
        "K7" = 549
    }
    'Irving Park' = @{
        "K3" = 43
        "K7" = 158
    }
    'River North' = @{
        "K3" = 48
        "K7" = 258
    }
    'Humboldt Park' = @{
        "K2" = 183
        "K3" = 214
        "K7" = 693
    }
    'Rogers Park' = @{
        "K2" = 64
        "K3" = 50
        "K6" = 77
        "K7" = 217
    }
    'Dunning' = @{
        "K4" = 8
        "K7" = 55
    }
    'Jefferson Park' = @{
        "K3" = 11
        "K7" = 39
    }
    'Washington Park' = @{
        "K5" = 4
        "K7" = 208
    }
    'Roseland' = @{
        "K2" = 158
        "K3" = 152
        "K7" = 471
    }
    'Near South Side' = @{
        "K4" = 15
        "K7" = 122
    }
    'Chicago Lawn' = @{
        "K6" = 126
        "K7" = 434
    }
    'Grand Boulevard' = @{
        "K2" = 97
        "K7" = 247
    }
    'West Loop' = @{
        "K3" = 48
        "K6" = 112
        "K7" = 253
    }
    'Kenwood' = @{
        "K2" = 37
        "K7" = 128
    }
    'Wicker Park' = @{
        "K3" = 18
        "K7" = 87
    }
    'Belmont Cragin' = @{
        "K2" = 120
        "K3" = 94
        "K7" = 356
    }
    'Gold Coast' = @{
        "K6" = 18
        "K7" = 30
    }
    'Portage Park' = @{
        "K3" = 41
        "K7" = 163
    }
    'Montclare' = @{
        "K4" = 4
        "K7" = 31
    }
    'Albany Park' = @{
        "K6" = 52
        "K7" = 164
    }
    'West Elsdon' = @{
        "K3" = 13
        "K7" = 69
    }
    'United Center' = @{
        "K2" = 54
        "K4" = 5
        "K7" = 206
    }
    'Uptown' = @{
        "K3" = 83
        "K7" = 272
    }
    'Edgewater' = @{
        "K2" = 49
        "K7" = 176
    }
    'Streeterville' = @{
        "K3" = 21
        "K7" = 126
    }
    'Washington Heights' = @{
        "K2" = 60
        "K4" = 12
        "K6" = 40
        "K7" = 170
    }
    'Hyde Park' = @{
        "K2" = 33
        "K7" = 166
    }
    'South Shore' = @{
        "K3" = 301
        "K6" = 215
        "K7" = 885
    }
    'Little Village' = @{
        "K6" = 179
        "K7" = 490
    }
    'Archer Heights' = @{
        "K2" = 25
        "K7" = 71
    }
    'Ukrainian Village' = @{
        "K6" = 17
        "K7" = 31
    }
}

foreach ($sheetName in $edits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $edits[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
